$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly data row (01-08-2021). Force column A to be entered as text
# so Excel doesn't auto-convert the "dd-mm-yyyy"-looking string into a
# date serial number; ClearFormats afterwards drops the temporary "@"
# number format so the cell keeps the workbook's default (unstyled) look,
# matching the rest of the data rows.
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "01-08-2021"
$ws.Range("A33").ClearFormats()

$ws.Range("B33").Value = 4.8
$ws.Range("C33").Value = 3.8
$ws.Range("D33").Value = 5.9
$ws.Range("E33").Value = 3.5
$ws.Range("F33").Value = 6.3
